$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Acknowledgments")

# Insert a new row at row 33 (the table is sorted alphabetically by column A,
# and "SudachiPy" sorts between "SSG"/"spaCy" row 32 and "TextBlob" row 33).
$ws.Rows.Item(33).Insert()

# Copy formatting from the row below (the old row 33, now row 34) so the new
# row matches the rest of the table's look (styles, number formats, etc.).
$ws.Rows.Item(34).Copy()
$ws.Rows.Item(33).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new row's values.
$ws.Cells.Item(33, 1).Value = "SudachiPy"
$ws.Cells.Item(33, 2).Value = "https://github.com/WorksApplications/sudachi.rs"
$ws.Cells.Item(33, 3).Value = "0.6.2"
$ws.Cells.Item(33, 4).Value = "Works Applications Co., Ltd."
$ws.Cells.Item(33, 5).Value = "Apache-2.0"
$ws.Cells.Item(33, 6).Value = "https://github.com/WorksApplications/sudachi.rs/blob/develop/LICENSE"

# Hyperlinks on the homepage (B) and license URL (F) cells.
$ws.Hyperlinks.Add($ws.Cells.Item(33, 2), "https://github.com/WorksApplications/sudachi.rs") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(33, 6), "https://github.com/WorksApplications/sudachi.rs/blob/develop/LICENSE") | Out-Null

# Restore the sorted range to include the new row and re-sort (keeps table
# consistent and matches the refreshed sortState range in the XML).
$ws.Sort.SortFields.Clear()
$sortRange = $ws.Range("A2:F37")
$keyRange = $ws.Range("A2:A37")
$ws.Sort.SortFields.Add($keyRange) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Update selection/view state to match the final saved state.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("C33").Select()
